# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-10-09 (serial 45208) to 2023-10-13 (serial 45212).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45208) {
        $cell.Value = 45212
    }
}
